# Remove the two bullet paragraphs "Suitable for small adenocarcinoma" and
# "Suitable for larger GI Stromal Tumors" from the "5 Distal Gastrectomy"
# section, leaving "Does not remove all lymph nodes" followed directly by
# the "Locally-advanced cancers..." paragraph.

$d = $word.ActiveDocument

# Locate the start of the first paragraph to remove.
$startRange = $d.Content
[void]$startRange.Find.Execute("Suitable for small adenocarcinoma", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Locate the end of the last paragraph to remove, then extend by one
# character so the deletion also swallows its trailing paragraph mark.
$endRange = $d.Content
[void]$endRange.Find.Execute("Suitable for larger GI Stromal Tumors", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$endRange.MoveEnd(1, 1)

# Build a range spanning both paragraphs (including paragraph marks) and
# delete it in one shot.
$deleteRange = $d.Range($startRange.Start, $endRange.End)
$deleteRange.Delete()
